## Updates the cryptos list (prices + 1h volume change %) to the latest
## scrape values, and swaps in new coins for the bottom of the table
## (THORChain/Aave reorder, EnergySwap -> FraxShare/ordi/MultiversX shift).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.916.97"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "'2.229.72"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'303.76"
$ws.Range("E5").Value = "  -4.43%  "
$ws.Range("D6").Value = "'93.58"
$ws.Range("E6").Value = "  -7.64%  "
$ws.Range("E7").Value = "  -1.96%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  -7.00%  "
$ws.Range("D10").Value = "'33.95"
$ws.Range("E10").Value = "  -8.17%  "
$ws.Range("E11").Value = "  -4.27%  "
$ws.Range("D12").Value = "'7.06"
$ws.Range("E12").Value = "  -8.06%  "
$ws.Range("E13").Value = "  -3.32%  "
$ws.Range("D14").Value = "'2.570.41"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "'2.260.54"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "'0.807"
$ws.Range("E16").Value = "  -5.72%  "
$ws.Range("D17").Value = "'13.43"
$ws.Range("E17").Value = "  -5.29%  "
$ws.Range("D18").Value = "'43.710.11"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  -3.61%  "
$ws.Range("D20").Value = "'12.09"
$ws.Range("E20").Value = "  -9.94%  "
$ws.Range("D21").Value = "'6.10"
$ws.Range("E21").Value = "  -6.59%  "
$ws.Range("D22").Value = "'64.08"
$ws.Range("E22").Value = "  -2.40%  "
$ws.Range("D23").Value = "'234.16"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  -7.44%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -10.24%  "
$ws.Range("D27").Value = "'9.71"
$ws.Range("E27").Value = "  -3.93%  "
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("D29").Value = "'35.63"
$ws.Range("E29").Value = "  -3.88%  "
$ws.Range("D30").Value = "'5.88"
$ws.Range("E30").Value = "  -6.01%  "
$ws.Range("D31").Value = "'19.67"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").Value = "'150.33"
$ws.Range("E32").Value = "  -5.03%  "
$ws.Range("D33").Value = "'0.0793"
$ws.Range("E33").Value = "  -6.68%  "
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("D35").Value = "'3.19"
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("D37").Value = "'0.106"
$ws.Range("E37").Value = "  -9.19%  "
$ws.Range("E38").Value = "  -11.46%  "
$ws.Range("D39").Value = "'14.41"
$ws.Range("E39").Value = "  -8.82%  "
$ws.Range("E40").Value = "  -10.59%  "
$ws.Range("E41").Value = "  -6.25%  "
$ws.Range("D42").Value = "'3.23"
$ws.Range("E42").Value = "  -13.38%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").Value = "'1.736.39"
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("D45").Value = "'83.25"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("E46").Value = "  -7.34%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'4.91"
$ws.Range("E47").Value = "  -5.48%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'98.91"
$ws.Range("E48").Value = "  -4.43%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'8.02"
$ws.Range("E49").Value = "  -4.13%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "'67.24"
$ws.Range("E50").Value = "  -10.94%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'53.43"
$ws.Range("E51").Value = "  -8.49%  "
